# Updates symbol list data (prices, swapped rows, labels) per the
# "Updated symbol list" GitHub Actions commit.
#
# NOTE: Column D (Price) is stored as text in this workbook (inline
# strings). Assigning a plain numeric-looking string via .Value would
# make Excel coerce it into a real number (dropping e.g. trailing
# zeros such as "0.05960" -> 0.0596), so every Price update below is
# entered with a leading single quote to force a text entry, exactly
# like typing '0.05960 into the cell through the Excel UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'251.48"
$ws.Range("D3").Value = "'23.84"
$ws.Range("D4").Value = "'5.945"
$ws.Range("D5").Value = "'0.05960"

# Row 6 / Row 7 swap (GateToken <-> KuCoinToken)
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D6").Value = "'6.570"
$ws.Range("E6").Value = "5KuCoinTokenKCS"

$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'3.416"
$ws.Range("E7").Value = "6GateTokenGT"

$ws.Range("D8").Value = "'1.321"
$ws.Range("D9").Value = "'0.7927"

# Row 10 / Row 11 / Row 12 cascade
# (One -> WazirX -> MandalaExchangeToken -> LiechtensteinCryptoassetsExchange)
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1482"
$ws.Range("E10").Value = "9WazirXWRX"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07839"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03355"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03033"
$ws.Range("E13").Value = "12BitrueCoinBTR"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09277"
$ws.Range("E14").Value = "13BitMartTokenBMX"

$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'3.558"
$ws.Range("E15").Value = "14MCDexMCB"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001668"
$ws.Range("E16").Value = "15BitForexTokenBF"

$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04779"
$ws.Range("E17").Value = "16CoinExTokenCET"

$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0006062"
$ws.Range("E18").Value = "17OneONE"

$ws.Range("D19").Value = "'0.006204"
$ws.Range("D20").Value = "'0.005687"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("D23").Value = "'3.689"
$ws.Range("D27").Value = "'0.0006477"
$ws.Range("D40").Value = "'0.04447"
$ws.Range("D41").Value = "'0.007019"

# Row 42 / Row 43 swap (CEJI <-> BKEXToken)
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1067"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003150"
$ws.Range("E43").Value = "42CEJICEJI"

$ws.Range("D44").Value = "'0.01041"

$ws.Range("E45").Value = "44ACDXExchangeACXTBestin24h"

$ws.Range("D46").Value = "'0.00005895"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D48").Value = "'0.7855"
$ws.Range("D49").Value = "'0.09847"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("D51").Value = "'0.01010"
